$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts the old row 13 "Segunda" data down to row 14)
$ws.Rows.Item(13).Insert()

# Row 12 becomes the new weekly data point (date moves forward a year, volume/prices/origin updated)
$ws.Cells.Item(12, 4).Value = 44627
$ws.Cells.Item(12, 13).Value = 16
$ws.Cells.Item(12, 14).Value = 405000
$ws.Cells.Item(12, 15).Value = 410000
$ws.Cells.Item(12, 16).Value = 407500
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 906

# The newly inserted row 13 re-creates the prior week's "Primera" grade record
# (identical to what row 12 held before this week's update)
$ws.Cells.Item(13, 1).Value = 8
$ws.Cells.Item(13, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44294
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100104
$ws.Cells.Item(13, 8).Value = "Frutos de pepita"
$ws.Cells.Item(13, 9).Value = 100104003
$ws.Cells.Item(13, 10).Value = "Membrillo"
$ws.Cells.Item(13, 11).Value = "Champion"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 20
$ws.Cells.Item(13, 14).Value = 225000
$ws.Cells.Item(13, 15).Value = 230000
$ws.Cells.Item(13, 16).Value = 227500
$ws.Cells.Item(13, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(13, 18).Value = "Región Metropolitana"
$ws.Cells.Item(13, 19).Value = 506
$ws.Cells.Item(13, 20).Value = 450
